# Commit: "Add files via upload"
# Team sheet: Prateek Tyagi's and Sudhansh Aggarwal's contact info was
# corrected/updated (new Stevens email + GitHub handle for Prateek; new
# Stevens email + GitHub handle for Sudhansh). Bob/Robert Majdi's row (RBM)
# is untouched.
$wb = $excel.ActiveWorkbook

$team = $wb.Worksheets.Item("Team")
$team.Range("B3").Value = "Prateek "
$team.Range("D3").Value = "ptyagi1@stevens.edu"
$team.Range("E3").Value = "pratt23"

$team.Range("D4").Value = "saggarw2@stevens.edu"
$team.Range("E4").Value = "agentweirdo"

# New email addresses get mailto hyperlinks, matching the existing one on D5.
$team.Hyperlinks.Add($team.Range("D3"), "mailto:ptyagi1@stevens.edu")
$team.Hyperlinks.Add($team.Range("D4"), "mailto:saggarw2@stevens.edu")

$team.Range("D15").Select() | Out-Null

# Sprint1 sheet: "Review Results" section edited — Keep-doing note changed
# from "bring alchohol; " to "R&D", and the Avoid note "Jim rowland " removed.
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Range("B18").Value = "R&D"
$sprint1.Range("B22").ClearContents()
$sprint1.Range("B18").Select() | Out-Null
